$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Second signature table ("${trabajador}" / "${vendedor}" row): the first
#    cell of the underline row becomes the "${trabajador_linea}" merge field,
#    and the now-orphaned "_GoBack" bookmark (Word always keeps one, marking
#    the last edit point) moves to sit right after that run.
# ---------------------------------------------------------------------------
$tabla = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidata = $d.Tables.Item($i)
    if ($candidata.Cell(2, 1).Range.Text -like "`${trabajador}*") {
        $tabla = $candidata
        break
    }
}
$celda = $tabla.Cell(1, 1)
$rCelda = $celda.Range
$rTexto = $d.Range($rCelda.Start, $rCelda.End - 1)

# Append a throw-away trailing character so the insertion point for the
# bookmark is not the very last position in the paragraph while we add it
# (avoids landing the bookmark at the wrong spot), then trim it back off.
$rTexto.Text = "`${trabajador_linea}X"
$posMarca = $rTexto.End - 1
$rMarca = $d.Range($posMarca, $posMarca)
$d.Bookmarks.Add("_GoBack", $rMarca)
$d.Range($posMarca, $posMarca + 1).Delete()

# ---------------------------------------------------------------------------
# 2) Drop the redundant "DNI: " label before the trabajador DNI merge field
#    (only the first occurrence in the document - the vendedor cell keeps
#    its own "DNI: " label untouched).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("DNI: `${", $true, $false, $false, $false, $false, `
    $true, 1, $false, "`${", 1) | Out-Null

# ---------------------------------------------------------------------------
# 3) A document only ever carries one "_GoBack" bookmark - re-adding it above
#    (step 1) already moved it off the trailing empty paragraph at the end
#    of the document, so that paragraph is bookmark-free again on its own;
#    nothing further to do there.
# ---------------------------------------------------------------------------
